$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 9 (before the roundLED group),
# shifting roundLED/rst/seg1..8 down by one row.
$ws.Rows.Item(9).Insert()

# Rebuild rows 3-69 with the final pin-table contents
# (diffLED/digitButtons/roundLED/seg1..8 groups reversed in row order,
# plus the new resetLED row).
$data = @{
    1 = @('clk', 'Input', 'PIN_Y2')
    2 = @('confirmButton', 'Input', 'PIN_R24')
    3 = @('diffLED[2]', 'Output', 'PIN_E25')
    4 = @('diffLED[1]', 'Output', 'PIN_E22')
    5 = @('diffLED[0]', 'Output', 'PIN_E21')
    6 = @('digitButtons[2]', 'Input', 'PIN_N21')
    7 = @('digitButtons[1]', 'Input', 'PIN_M21')
    8 = @('digitButtons[0]', 'Input', 'PIN_M23')
    9 = @('resetLED', 'Output', 'PIN_G19')
    10 = @('roundLED[2]', 'Output', 'PIN_H15')
    11 = @('roundLED[1]', 'Output', 'PIN_G16')
    12 = @('roundLED[0]', 'Output', 'PIN_G15')
    13 = @('rst', 'Input', 'PIN_AB28')
    14 = @('seg1[6]', 'Output', 'PIN_AA14')
    15 = @('seg1[5]', 'Output', 'PIN_AG18')
    16 = @('seg1[4]', 'Output', 'PIN_AF17')
    17 = @('seg1[3]', 'Output', 'PIN_AH17')
    18 = @('seg1[2]', 'Output', 'PIN_AG17')
    19 = @('seg1[1]', 'Output', 'PIN_AE17')
    20 = @('seg1[0]', 'Output', 'PIN_AD17')
    21 = @('seg2[6]', 'Output', 'PIN_AC17')
    22 = @('seg2[5]', 'Output', 'PIN_AA15')
    23 = @('seg2[4]', 'Output', 'PIN_AB15')
    24 = @('seg2[3]', 'Output', 'PIN_AB17')
    25 = @('seg2[2]', 'Output', 'PIN_AA16')
    26 = @('seg2[1]', 'Output', 'PIN_AB16')
    27 = @('seg2[0]', 'Output', 'PIN_AA17')
    28 = @('seg3[6]', 'Output', 'PIN_AH18')
    29 = @('seg3[5]', 'Output', 'PIN_AF18')
    30 = @('seg3[4]', 'Output', 'PIN_AG19')
    31 = @('seg3[3]', 'Output', 'PIN_AH19')
    32 = @('seg3[2]', 'Output', 'PIN_AB18')
    33 = @('seg3[1]', 'Output', 'PIN_AC18')
    34 = @('seg3[0]', 'Output', 'PIN_AD18')
    35 = @('seg4[6]', 'Output', 'PIN_AE18')
    36 = @('seg4[5]', 'Output', 'PIN_AF19')
    37 = @('seg4[4]', 'Output', 'PIN_AE19')
    38 = @('seg4[3]', 'Output', 'PIN_AH21')
    39 = @('seg4[2]', 'Output', 'PIN_AG21')
    40 = @('seg4[1]', 'Output', 'PIN_AA19')
    41 = @('seg4[0]', 'Output', 'PIN_AB19')
    42 = @('seg5[6]', 'Output', 'PIN_Y19')
    43 = @('seg5[5]', 'Output', 'PIN_AF23')
    44 = @('seg5[4]', 'Output', 'PIN_AD24')
    45 = @('seg5[3]', 'Output', 'PIN_AA21')
    46 = @('seg5[2]', 'Output', 'PIN_AB20')
    47 = @('seg5[1]', 'Output', 'PIN_U21')
    48 = @('seg5[0]', 'Output', 'PIN_V21')
    49 = @('seg6[6]', 'Output', 'PIN_W28')
    50 = @('seg6[5]', 'Output', 'PIN_W27')
    51 = @('seg6[4]', 'Output', 'PIN_Y26')
    52 = @('seg6[3]', 'Output', 'PIN_W26')
    53 = @('seg6[2]', 'Output', 'PIN_Y25')
    54 = @('seg6[1]', 'Output', 'PIN_AA26')
    55 = @('seg6[0]', 'Output', 'PIN_AA25')
    56 = @('seg7[6]', 'Output', 'PIN_U24')
    57 = @('seg7[5]', 'Output', 'PIN_U23')
    58 = @('seg7[4]', 'Output', 'PIN_W25')
    59 = @('seg7[3]', 'Output', 'PIN_W22')
    60 = @('seg7[2]', 'Output', 'PIN_W21')
    61 = @('seg7[1]', 'Output', 'PIN_Y22')
    62 = @('seg7[0]', 'Output', 'PIN_M24')
    63 = @('seg8[6]', 'Output', 'PIN_H22')
    64 = @('seg8[5]', 'Output', 'PIN_J22')
    65 = @('seg8[4]', 'Output', 'PIN_L25')
    66 = @('seg8[3]', 'Output', 'PIN_L26')
    67 = @('seg8[2]', 'Output', 'PIN_E17')
    68 = @('seg8[1]', 'Output', 'PIN_F22')
    69 = @('seg8[0]', 'Output', 'PIN_G18')
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item([int]$r, 1).Value = $vals[0]
    $ws.Cells.Item([int]$r, 2).Value = $vals[1]
    $ws.Cells.Item([int]$r, 3).Value = $vals[2]
}

# Update the saved selection to E3 (matches the target view state)
$ws.Range("E3").Select()

Write-Host "done"
